$d = $word.ActiveDocument

$pairs = @(
    @("56×56=", "59×93="),
    @("31×82=", "95×27="),
    @("33×43=", "36×21="),
    @("17×32=", "51×35="),
    @("33×97=", "87×97="),
    @("27×49=", "16×67="),
    @("41×85=", "22×97="),
    @("99×85=", "22×70="),
    @("98×31=", "89×14="),
    @("18×84=", "58×42="),
    @("58×77=", "46×29="),
    @("71×69=", "20×63="),
    @("99×53=", "52×97="),
    @("19×41=", "16×46="),
    @("63×77=", "20×12="),
    @("43×50=", "76×48="),
    @("31×41=", "36×78="),
    @("56×35=", "56×40="),
    @("89×78=", "37×38="),
    @("92×92=", "53×18="),
    @("75×60=", "78×25="),
    @("71×19=", "35×39="),
    @("79×23=", "13×71="),
    @("99×88=", "37×58="),
    @("26×98=", "71×72=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
